# Remove the trailing "Ver no Jupiter..." / "(c) 2020 ..." footer block
# (plus the blank paragraph immediately preceding it), leaving the blank
# paragraph that used to follow the footer intact.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph by
# scanning paragraph text (robust against any renumbering).
$jupiterIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Ver no Jupiter*") {
        $jupiterIdx = $i
        break
    }
}

if ($jupiterIdx -eq -1) {
    throw "Could not locate the 'Ver no Jupiter' paragraph"
}

# Delete the blank paragraph right before it, the "Ver no Jupiter..."
# paragraph itself, and the "(c) 2020 ..." paragraph right after it.
$firstToRemove = $d.Paragraphs.Item($jupiterIdx - 1)
$lastToRemove = $d.Paragraphs.Item($jupiterIdx + 1)

$deleteRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
$deleteRange.Delete()
